$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the header row (row 1) with two new columns inserted after "LandId"
# (Type, LandCover), the remaining original headers shifted right, the old
# "AreaMeasurementId" column dropped, and two new trailing columns added
# (GeorefId, Guid); "GUID" is effectively replaced by "Guid" at the end.
$ws.Range("A1").Value = "LandId"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "LandCover"
$ws.Range("D1").Value = "ParcelNr."
$ws.Range("E1").Value = "LandCode"
$ws.Range("F1").Value = "ValidFrom"
$ws.Range("G1").Value = "ValidUntil"
$ws.Range("H1").Value = "Name"
$ws.Range("I1").Value = "SiteId"
$ws.Range("J1").Value = "GeorefId"
$ws.Range("K1").Value = "Guid"

# Give the header row its own (applied-font) style.
$ws.Range("A1:K1").Font.Bold = $true

# Update the active selection to E6, matching where the user clicked next.
$ws.Range("E6").Select()
